$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260590672492981
$ws.Range("B1").Value = 3.857060432434082
$ws.Range("C1").Value = 3.449803113937378
$ws.Range("D1").Value = 3.702497005462646
$ws.Range("E1").Value = 1.189799904823303
